$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (2-10) for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$data = @{
  2  = @{ E=3; G=11.779764; H=35.339292; I=0.1028447940505417; J=0.1028447940505417; K=3; M=59.30472066666667; N=177.914162; O=0.3117418698225773; P=0.3117418698225772; Q=698.5956135392561; R=6287.360521853304; S=0.03206102839883375; T=0.03206102839883374 }
  3  = @{ E=3; G=11.779764; H=35.339292; I=0.1028447940505417; J=0.1028447940505417; K=3; M=99.97060400000002; N=299.9118120000001; O=0.5255066151212701; P=0.5255066151212699; Q=1177.630122057456; R=10598.67109851711; S=0.05404561960434431; T=0.05404561960434429 }
  4  = @{ E=3; G=11.779764; H=35.339292; I=0.1028447940505417; J=0.1028447940505417; K=3; M=30.96129866666667; N=92.88389600000001; O=0.1627515150561528; P=0.1627515150561528; Q=364.716791426848; R=3282.451122841632; S=0.01673814604736367; T=0.01673814604736367 }
  5  = @{ E=3; G=43.841352; H=131.524056; I=0.3827627461243965; J=0.3827627461243964; K=3; M=59.30472066666667; N=177.914162; O=0.3117418698225773; P=0.3117418698225772; Q=2599.999134009008; R=23399.99220608107; S=0.1193231741752438; T=0.1193231741752438 }
  6  = @{ E=3; G=43.841352; H=131.524056; I=0.3827627461243965; J=0.3827627461243964; K=3; M=99.97060400000002; N=299.9118120000001; O=0.5255066151212701; P=0.5255066151212699; Q=4382.846439616609; R=39445.61795654948; S=0.2011443551103536; T=0.2011443551103536 }
  7  = @{ E=3; G=43.841352; H=131.524056; I=0.3827627461243965; J=0.3827627461243964; K=3; M=30.96129866666667; N=92.88389600000001; O=0.1627515150561528; P=0.1627515150561528; Q=1357.385193222464; R=12216.46673900218; S=0.0622952168387991; T=0.06229521683879908 }
  8  = @{ E=3; G=58.91811866666666; H=176.754356; I=0.5143924598250619; J=0.5143924598250619; K=3; M=59.30472066666667; N=177.914162; O=0.3117418698225773; P=0.3117418698225772; Q=3494.122569732186; R=31447.10312758967; S=0.1603576672484998; T=0.1603576672484997 }
  9  = @{ E=3; G=58.91811866666666; H=176.754356; I=0.5143924598250619; J=0.5143924598250619; K=3; M=99.97060400000002; N=299.9118120000001; O=0.5255066151212701; P=0.5255066151212699; Q=5890.079909650342; R=53010.71918685308; S=0.2703166404065722; T=0.2703166404065721 }
  10 = @{ E=3; G=58.91811866666666; H=176.754356; I=0.5143924598250619; J=0.5143924598250619; K=3; M=30.96129866666667; N=92.88389600000001; O=0.1627515150561528; P=0.1627515150561528; Q=1824.181468916775; R=16417.63322025098; S=0.08371815216999003; T=0.08371815216999001 }
}

foreach ($row in $data.Keys) {
  $cols = $data[$row]
  foreach ($col in $cols.Keys) {
    $ws.Range("$col$row").Value = $cols[$col]
  }
}
